$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. document.xml: the sole (empty) paragraph loses its NoIndent
#    paragraph style, reverting to a plain default paragraph.
# ------------------------------------------------------------------
$p = $d.Paragraphs(1)
$p.Style = $d.Styles("Normal")

# ------------------------------------------------------------------
# 2. styles.xml: link the existing "Block Text" style to a new
#    "Block Text Char" character style.
# ------------------------------------------------------------------
$blockText = $d.Styles("Block Text")
$blockText.LinkStyle = "BlockTextChar"

# ------------------------------------------------------------------
# 3. styles.xml: add the new paragraph style "NextBlockText", based
#    on BlockText, linked to a new "NextBlockText Char" style, with a
#    first-line indent of 720 twips (36 pt) so paragraphs after the
#    first one in a block quote get indented.
# ------------------------------------------------------------------
$nextBlockText = $d.Styles.Add("NextBlockText", 1)
$nextBlockText.NameLocal = "NextBlockText"
$nextBlockText.BaseStyle = "BlockText"
$nextBlockText.LinkStyle = "NextBlockTextChar"
$nextBlockText.QuickStyle = $true
$nextBlockText.ParagraphFormat.FirstLineIndent = 36

# ------------------------------------------------------------------
# 4. styles.xml: add the new character style "BlockTextChar" (linked
#    back to BlockText), based on BodyTextChar.
# ------------------------------------------------------------------
$blockTextChar = $d.Styles.Add("BlockTextChar", 2)
$blockTextChar.NameLocal = "Block Text Char"
$blockTextChar.BaseStyle = "BodyTextChar"
$blockTextChar.LinkStyle = "BlockText"
$blockTextChar.Priority = 9

# ------------------------------------------------------------------
# 5. styles.xml: add the new character style "NextBlockTextChar",
#    based on BlockTextChar, linked back to NextBlockText.
# ------------------------------------------------------------------
$nextBlockTextChar = $d.Styles.Add("NextBlockTextChar", 2)
$nextBlockTextChar.NameLocal = "NextBlockText Char"
$nextBlockTextChar.BaseStyle = "BlockTextChar"
$nextBlockTextChar.LinkStyle = "NextBlockText"
